# Generate Report for Handoff
# - Flip the "Status" from "In Translation" to "Ready for handoff"
#   (Overview!E2/F2, zh-cn!C2, de-de!C2).
# - Refresh the associated handoff timestamps by ~30s
#   (Overview!G2, de-de!H2 share one string; zh-cn!H2 is separate).
# - Widen the now-longer "Status"/language summary columns so the new
#   text isn't clipped (Overview E:F, zh-cn C, de-de C).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Timestamps ---
$wsOverview.Range("G2").Value = "2016-08-18 00:54:49"
$wsDeDe.Range("H2").Value     = "2016-08-18 00:54:49"
$wsZhCn.Range("H2").Value     = "2016-08-18 00:54:45"

# --- Column widths: widen to fit "Ready for handoff" ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3826548258464
$wsOverview.Columns.Item(6).ColumnWidth = 16.3826548258464
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.3826548258464
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.3826548258464
